# Generate Report for Handback
# This script re-runs the localization status report generation: the
# 4e26e95f-30d7-4575-b028-0c8920bc4c7f.md file has now been handed back
# (in sync with en-US), swapping positions/content with
# f9ecdccb-1235-4204-82ea-06d7c584e065.md across the Overview, zh-cn and
# de-de sheets, and refreshing the relevant handback timestamps.

function Set-HyperlinkText($ws, $addr, $text) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.TextToDisplay = $text
        }
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A4").Value2 = "4e26e95f-30d7-4575-b028-0c8920bc4c7f.md"
$ws1.Range("B4").Value2 = "e2e\4e26e95f-30d7-4575-b028-0c8920bc4c7f.md"
$ws1.Range("G4").Value2 = "2016-10-13 13:25:08"

$ws1.Range("A5").Value2 = "f9ecdccb-1235-4204-82ea-06d7c584e065.md"
$ws1.Range("B5").Value2 = "e2e\f9ecdccb-1235-4204-82ea-06d7c584e065.md"
$ws1.Range("E5").Value2 = "Handed back: in sync with en-US"
$ws1.Range("F5").Value2 = "Handed back: in sync with en-US"
$ws1.Range("G5").Value2 = "2016-10-13 13:22:48"

Set-HyperlinkText $ws1 "`$B`$4" "e2e\4e26e95f-30d7-4575-b028-0c8920bc4c7f.md"
Set-HyperlinkText $ws1 "`$B`$5" "e2e\f9ecdccb-1235-4204-82ea-06d7c584e065.md"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A4").Value2 = "4e26e95f-30d7-4575-b028-0c8920bc4c7f.md"
$ws2.Range("F4").Value2 = "'False"
$ws2.Range("G4").Value2 = "4e26e95f-30d7-4575-b028-0c8920bc4c7f.aae293ece568c5c5c6773c52687e73db68f8c4f3.zh-cn.xlf"
$ws2.Range("H4").Value2 = "2016-10-13 13:24:57"
$ws2.Range("I4").Value2 = "4e26e95f-30d7-4575-b028-0c8920bc4c7f.md"
$ws2.Range("J4").Value2 = "4e26e95f-30d7-4575-b028-0c8920bc4c7f.aae293ece568c5c5c6773c52687e73db68f8c4f3.zh-cn.xlf"
$ws2.Range("K4").Value2 = "2016-10-13 13:26:34"

$ws2.Range("A5").Value2 = "f9ecdccb-1235-4204-82ea-06d7c584e065.md"
$ws2.Range("C5").Value2 = "Handed back: in sync with en-US"
$ws2.Range("F5").Value2 = "'True"
$ws2.Range("G5").Value2 = "2a13abd4-60e7-4874-99f6-09924cf324fa.7e32b9a3179c819fe8848520ce3d2106f6a5dc7c.zh-cn.xlf"
$ws2.Range("H5").Value2 = "2016-10-13 13:22:37"
$ws2.Range("I5").Value2 = "2a13abd4-60e7-4874-99f6-09924cf324fa.md"
$ws2.Range("J5").Value2 = "2a13abd4-60e7-4874-99f6-09924cf324fa.7e32b9a3179c819fe8848520ce3d2106f6a5dc7c.zh-cn.xlf"
$ws2.Range("K5").Value2 = "2016-10-13 13:23:18"
$ws2.Range("P5").Value2 = ""

Set-HyperlinkText $ws2 "`$A`$4" "4e26e95f-30d7-4575-b028-0c8920bc4c7f.md"
Set-HyperlinkText $ws2 "`$I`$4" "4e26e95f-30d7-4575-b028-0c8920bc4c7f.md"
Set-HyperlinkText $ws2 "`$A`$5" "f9ecdccb-1235-4204-82ea-06d7c584e065.md"
Set-HyperlinkText $ws2 "`$I`$5" "2a13abd4-60e7-4874-99f6-09924cf324fa.md"

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A4").Value2 = "4e26e95f-30d7-4575-b028-0c8920bc4c7f.md"
$ws3.Range("F4").Value2 = "'False"
$ws3.Range("G4").Value2 = "4e26e95f-30d7-4575-b028-0c8920bc4c7f.aae293ece568c5c5c6773c52687e73db68f8c4f3.de-de.xlf"
$ws3.Range("H4").Value2 = "2016-10-13 13:25:08"
$ws3.Range("I4").Value2 = "4e26e95f-30d7-4575-b028-0c8920bc4c7f.md"
$ws3.Range("J4").Value2 = "4e26e95f-30d7-4575-b028-0c8920bc4c7f.aae293ece568c5c5c6773c52687e73db68f8c4f3.de-de.xlf"
$ws3.Range("K4").Value2 = "2016-10-13 13:26:50"

$ws3.Range("A5").Value2 = "f9ecdccb-1235-4204-82ea-06d7c584e065.md"
$ws3.Range("C5").Value2 = "Handed back: in sync with en-US"
$ws3.Range("F5").Value2 = "'True"
$ws3.Range("G5").Value2 = "2a13abd4-60e7-4874-99f6-09924cf324fa.7e32b9a3179c819fe8848520ce3d2106f6a5dc7c.de-de.xlf"
$ws3.Range("H5").Value2 = "2016-10-13 13:22:48"
$ws3.Range("I5").Value2 = "2a13abd4-60e7-4874-99f6-09924cf324fa.md"
$ws3.Range("J5").Value2 = "2a13abd4-60e7-4874-99f6-09924cf324fa.7e32b9a3179c819fe8848520ce3d2106f6a5dc7c.de-de.xlf"
$ws3.Range("K5").Value2 = "2016-10-13 13:23:35"
$ws3.Range("P5").Value2 = ""

Set-HyperlinkText $ws3 "`$A`$4" "4e26e95f-30d7-4575-b028-0c8920bc4c7f.md"
Set-HyperlinkText $ws3 "`$I`$4" "4e26e95f-30d7-4575-b028-0c8920bc4c7f.md"
Set-HyperlinkText $ws3 "`$A`$5" "f9ecdccb-1235-4204-82ea-06d7c584e065.md"
Set-HyperlinkText $ws3 "`$I`$5" "2a13abd4-60e7-4874-99f6-09924cf324fa.md"

# ---------------------------------------------------------------
# Column P width shrinks back to default now that the long error
# message text is gone from zh-cn / de-de sheets.
# ---------------------------------------------------------------
$ws2.Columns.Item(16).ColumnWidth = 13.7470531463623
$ws3.Columns.Item(16).ColumnWidth = 13.7470531463623

$wb.Save()
